# Add team record (Wins/Losses/Ties) columns to the HOU_2008 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 - copy formatting (bold, border, alignment) from
# the existing header style used in A1:AC1 (e.g. AC1), then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record for every player row (2-43): 86 wins, 75 losses, 0 ties.
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 86
    $ws.Cells.Item($r, 31).Value = 75
    $ws.Cells.Item($r, 32).Value = 0
}
